$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 235-239 (Romania Liga I, 2024-06-02 slate) ---
# Row 235
$ws.Cells.Item(235, 2).Value = 6870268
$ws.Cells.Item(235, 5).Value = "Petrolul Ploiesti"
$ws.Cells.Item(235, 6).Value = "ACS Sepsi"
$ws.Cells.Item(235, 7).Value = 1
$ws.Cells.Item(235, 8).Value = 2
$ws.Cells.Item(235, 9).Value = 0
$ws.Cells.Item(235, 10).Value = 1
$ws.Cells.Item(235, 11).Value = "A"
$ws.Cells.Item(235, 12).Value = 2.8
$ws.Cells.Item(235, 13).Value = 3
$ws.Cells.Item(235, 14).Value = 2.55
$ws.Cells.Item(235, 15).Value = 3
$ws.Cells.Item(235, 16).Value = 3.2
$ws.Cells.Item(235, 17).Value = 2.3
$ws.Cells.Item(235, 18).Value = 0.25
$ws.Cells.Item(235, 19).Value = 1.85
$ws.Cells.Item(235, 20).Value = 2
$ws.Cells.Item(235, 21).Value = 2.25
$ws.Cells.Item(235, 22).Value = 1.875
$ws.Cells.Item(235, 23).Value = 1.975
$ws.Cells.Item(235, 24).Value = -1
$ws.Cells.Item(235, 25).Value = -1
$ws.Cells.Item(235, 26).Value = 1.3
$ws.Cells.Item(235, 27).Value = -1
$ws.Cells.Item(235, 28).Value = 1
$ws.Cells.Item(235, 29).Value = 0.875
$ws.Cells.Item(235, 30).Value = -1

# Row 236
$ws.Cells.Item(236, 2).Value = 6865915
$ws.Cells.Item(236, 5).Value = "FC Voluntari"
$ws.Cells.Item(236, 6).Value = "Universitatea Cluj"
$ws.Cells.Item(236, 7).Value = 0
$ws.Cells.Item(236, 8).Value = 0
$ws.Cells.Item(236, 9).Value = 0
$ws.Cells.Item(236, 10).Value = 0
$ws.Cells.Item(236, 11).Value = "D"
$ws.Cells.Item(236, 12).Value = 3.5
$ws.Cells.Item(236, 13).Value = 3.25
$ws.Cells.Item(236, 14).Value = 2.05
$ws.Cells.Item(236, 15).Value = 3.4
$ws.Cells.Item(236, 16).Value = 3.1
$ws.Cells.Item(236, 17).Value = 2.15
$ws.Cells.Item(236, 18).Value = 0.25
$ws.Cells.Item(236, 19).Value = 1.975
$ws.Cells.Item(236, 20).Value = 1.875
$ws.Cells.Item(236, 21).Value = 2.25
$ws.Cells.Item(236, 22).Value = 2.05
$ws.Cells.Item(236, 23).Value = 1.75
$ws.Cells.Item(236, 24).Value = -1
$ws.Cells.Item(236, 25).Value = 2.1
$ws.Cells.Item(236, 26).Value = -1
$ws.Cells.Item(236, 27).Value = 0.4875
$ws.Cells.Item(236, 28).Value = -0.5
$ws.Cells.Item(236, 29).Value = -1
$ws.Cells.Item(236, 30).Value = 0.75

# Row 237
$ws.Cells.Item(237, 2).Value = 6861095
$ws.Cells.Item(237, 5).Value = "FC Botosani"
$ws.Cells.Item(237, 6).Value = "Farul Constanta"
$ws.Cells.Item(237, 7).Value = 0
$ws.Cells.Item(237, 8).Value = 0
$ws.Cells.Item(237, 9).Value = 0
$ws.Cells.Item(237, 10).Value = 0
$ws.Cells.Item(237, 11).Value = "D"
$ws.Cells.Item(237, 12).Value = 3.75
$ws.Cells.Item(237, 13).Value = 3.4
$ws.Cells.Item(237, 14).Value = 1.909
$ws.Cells.Item(237, 15).Value = 3.1
$ws.Cells.Item(237, 16).Value = 3
$ws.Cells.Item(237, 17).Value = 2.375
$ws.Cells.Item(237, 18).Value = 0.25
$ws.Cells.Item(237, 19).Value = 1.775
$ws.Cells.Item(237, 20).Value = 2.1
$ws.Cells.Item(237, 21).Value = 2
$ws.Cells.Item(237, 22).Value = 1.8
$ws.Cells.Item(237, 23).Value = 2.05
$ws.Cells.Item(237, 24).Value = -1
$ws.Cells.Item(237, 25).Value = 2
$ws.Cells.Item(237, 26).Value = -1
$ws.Cells.Item(237, 27).Value = 0.3875
$ws.Cells.Item(237, 28).Value = -0.5
$ws.Cells.Item(237, 29).Value = -1
$ws.Cells.Item(237, 30).Value = 1.05

# Row 238
$ws.Cells.Item(238, 2).Value = 6836277
$ws.Cells.Item(238, 5).Value = "CFR Cluj"
$ws.Cells.Item(238, 6).Value = "AFC Hermannstadt"
$ws.Cells.Item(238, 7).Value = 1
$ws.Cells.Item(238, 8).Value = 0
$ws.Cells.Item(238, 9).Value = 0
$ws.Cells.Item(238, 10).Value = 0
$ws.Cells.Item(238, 11).Value = "H"
$ws.Cells.Item(238, 12).Value = 1.7
$ws.Cells.Item(238, 13).Value = 3.4
$ws.Cells.Item(238, 14).Value = 5
$ws.Cells.Item(238, 15).Value = 1.65
$ws.Cells.Item(238, 16).Value = 3.5
$ws.Cells.Item(238, 17).Value = 5.25
$ws.Cells.Item(238, 18).Value = -0.75
$ws.Cells.Item(238, 19).Value = 1.85
$ws.Cells.Item(238, 20).Value = 2
$ws.Cells.Item(238, 21).Value = 2.25
$ws.Cells.Item(238, 22).Value = 1.875
$ws.Cells.Item(238, 23).Value = 1.975
$ws.Cells.Item(238, 24).Value = 0.6499999999999999
$ws.Cells.Item(238, 25).Value = -1
$ws.Cells.Item(238, 26).Value = -1
$ws.Cells.Item(238, 27).Value = 0.425
$ws.Cells.Item(238, 28).Value = -0.5
$ws.Cells.Item(238, 29).Value = -1
$ws.Cells.Item(238, 30).Value = 0.9750000000000001

# Row 239
$ws.Cells.Item(239, 2).Value = 6852370
$ws.Cells.Item(239, 5).Value = "Dinamo Bucharest"
$ws.Cells.Item(239, 6).Value = "ACS UTA Batrana Doamna"
$ws.Cells.Item(239, 7).Value = 1
$ws.Cells.Item(239, 8).Value = 0
$ws.Cells.Item(239, 9).Value = 1
$ws.Cells.Item(239, 10).Value = 0
$ws.Cells.Item(239, 11).Value = "H"
$ws.Cells.Item(239, 12).Value = 2.55
$ws.Cells.Item(239, 13).Value = 2.875
$ws.Cells.Item(239, 14).Value = 3
$ws.Cells.Item(239, 15).Value = 2.375
$ws.Cells.Item(239, 16).Value = 3
$ws.Cells.Item(239, 17).Value = 3.1
$ws.Cells.Item(239, 18).Value = -0.25
$ws.Cells.Item(239, 19).Value = 2
$ws.Cells.Item(239, 20).Value = 1.85
$ws.Cells.Item(239, 21).Value = 2.25
$ws.Cells.Item(239, 22).Value = 1.975
$ws.Cells.Item(239, 23).Value = 1.875
$ws.Cells.Item(239, 24).Value = 1.375
$ws.Cells.Item(239, 25).Value = -1
$ws.Cells.Item(239, 26).Value = -1
$ws.Cells.Item(239, 27).Value = 1
$ws.Cells.Item(239, 28).Value = -1
$ws.Cells.Item(239, 29).Value = -1
$ws.Cells.Item(239, 30).Value = 0.875

# --- Rows 309-312 (Romania Liga I, 2024-06-03 slate) ---
# Row 309
$ws.Cells.Item(309, 2).Value = 8191523
$ws.Cells.Item(309, 5).Value = "Otelul Galati"
$ws.Cells.Item(309, 6).Value = "FC Botosani"
$ws.Cells.Item(309, 7).Value = 2
$ws.Cells.Item(309, 8).Value = 0
$ws.Cells.Item(309, 9).Value = 2
$ws.Cells.Item(309, 10).Value = 0
$ws.Cells.Item(309, 11).Value = "H"
$ws.Cells.Item(309, 12).Value = 1.666
$ws.Cells.Item(309, 13).Value = 3.6
$ws.Cells.Item(309, 14).Value = 4.6
$ws.Cells.Item(309, 15).Value = 2.9
$ws.Cells.Item(309, 16).Value = 3.5
$ws.Cells.Item(309, 17).Value = 2.2
$ws.Cells.Item(309, 18).Value = 0.25
$ws.Cells.Item(309, 19).Value = 1.85
$ws.Cells.Item(309, 20).Value = 2
$ws.Cells.Item(309, 21).Value = 2.25
$ws.Cells.Item(309, 22).Value = 1.875
$ws.Cells.Item(309, 23).Value = 1.975
$ws.Cells.Item(309, 24).Value = 1.9
$ws.Cells.Item(309, 25).Value = -1
$ws.Cells.Item(309, 26).Value = -1
$ws.Cells.Item(309, 27).Value = 0.8500000000000001
$ws.Cells.Item(309, 28).Value = -1
$ws.Cells.Item(309, 29).Value = -0.5
$ws.Cells.Item(309, 30).Value = 0.4875

# Row 310
$ws.Cells.Item(310, 2).Value = 8191463
$ws.Cells.Item(310, 5).Value = "Dinamo Bucharest"
$ws.Cells.Item(310, 6).Value = "ACS UTA Batrana Doamna"
$ws.Cells.Item(310, 7).Value = 2
$ws.Cells.Item(310, 8).Value = 0
$ws.Cells.Item(310, 9).Value = 2
$ws.Cells.Item(310, 10).Value = 0
$ws.Cells.Item(310, 11).Value = "H"
$ws.Cells.Item(310, 12).Value = 1.833
$ws.Cells.Item(310, 13).Value = 3.4
$ws.Cells.Item(310, 14).Value = 3.6
$ws.Cells.Item(310, 15).Value = 1.5
$ws.Cells.Item(310, 16).Value = 4.333
$ws.Cells.Item(310, 17).Value = 5
$ws.Cells.Item(310, 18).Value = -1
$ws.Cells.Item(310, 19).Value = 1.875
$ws.Cells.Item(310, 20).Value = 1.975
$ws.Cells.Item(310, 21).Value = 3
$ws.Cells.Item(310, 22).Value = 2.025
$ws.Cells.Item(310, 23).Value = 1.825
$ws.Cells.Item(310, 24).Value = 0.5
$ws.Cells.Item(310, 25).Value = -1
$ws.Cells.Item(310, 26).Value = -1
$ws.Cells.Item(310, 27).Value = 0.875
$ws.Cells.Item(310, 28).Value = -1
$ws.Cells.Item(310, 29).Value = -1
$ws.Cells.Item(310, 30).Value = 0.825

# Row 311
$ws.Cells.Item(311, 2).Value = 8191475
$ws.Cells.Item(311, 5).Value = "FC U Craiova 1948"
$ws.Cells.Item(311, 6).Value = "AFC Hermannstadt"
$ws.Cells.Item(311, 7).Value = 1
$ws.Cells.Item(311, 8).Value = 3
$ws.Cells.Item(311, 9).Value = 0
$ws.Cells.Item(311, 10).Value = 0
$ws.Cells.Item(311, 11).Value = "A"
$ws.Cells.Item(311, 12).Value = 2.625
$ws.Cells.Item(311, 13).Value = 3.3
$ws.Cells.Item(311, 14).Value = 2.45
$ws.Cells.Item(311, 15).Value = 2.05
$ws.Cells.Item(311, 16).Value = 3.5
$ws.Cells.Item(311, 17).Value = 3
$ws.Cells.Item(311, 18).Value = -0.25
$ws.Cells.Item(311, 19).Value = 1.85
$ws.Cells.Item(311, 20).Value = 2
$ws.Cells.Item(311, 21).Value = 2.25
$ws.Cells.Item(311, 22).Value = 1.825
$ws.Cells.Item(311, 23).Value = 2.025
$ws.Cells.Item(311, 24).Value = -1
$ws.Cells.Item(311, 25).Value = -1
$ws.Cells.Item(311, 26).Value = 2
$ws.Cells.Item(311, 27).Value = -1
$ws.Cells.Item(311, 28).Value = 1
$ws.Cells.Item(311, 29).Value = 0.825
$ws.Cells.Item(311, 30).Value = -1

# Row 312
$ws.Cells.Item(312, 2).Value = 8191462
$ws.Cells.Item(312, 5).Value = "CSM Politehnica Iasi"
$ws.Cells.Item(312, 6).Value = "Petrolul Ploiesti"
$ws.Cells.Item(312, 7).Value = 2
$ws.Cells.Item(312, 8).Value = 0
$ws.Cells.Item(312, 9).Value = 0
$ws.Cells.Item(312, 10).Value = 0
$ws.Cells.Item(312, 11).Value = "H"
$ws.Cells.Item(312, 12).Value = 2.1
$ws.Cells.Item(312, 13).Value = 3.3
$ws.Cells.Item(312, 14).Value = 3.1
$ws.Cells.Item(312, 15).Value = 1.8
$ws.Cells.Item(312, 16).Value = 3.2
$ws.Cells.Item(312, 17).Value = 4.2
$ws.Cells.Item(312, 18).Value = -0.5
$ws.Cells.Item(312, 19).Value = 1.85
$ws.Cells.Item(312, 20).Value = 2
$ws.Cells.Item(312, 21).Value = 2.25
$ws.Cells.Item(312, 22).Value = 2.025
$ws.Cells.Item(312, 23).Value = 1.825
$ws.Cells.Item(312, 24).Value = 0.8
$ws.Cells.Item(312, 25).Value = -1
$ws.Cells.Item(312, 26).Value = -1
$ws.Cells.Item(312, 27).Value = 0.8500000000000001
$ws.Cells.Item(312, 28).Value = -1
$ws.Cells.Item(312, 29).Value = -0.5
$ws.Cells.Item(312, 30).Value = 0.4125
